$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,4: Employment volumes / by occupation / by industry (Annual Population Survey)
# C column = "Latest period (release date)", D column = "Next period (release date)"
$ws.Range("C2:C4").Value = "Jul 2022 - Jun 2023 (24/10/23)"
$ws.Range("D2:D4").Value = "Oct 2022 - Sep 2023 (16/01/24)"

# Rows 8,9: Enterprises by employment size band / industry (ONS UK Business Counts)
$ws.Range("D8:D9").Value = "Mar 2024 (Sep 24)"

# Update the active selection to match the author's final cursor position
$ws.Range("D6").Select()
